$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 41
$ws.Range("E41").Value = 23
$ws.Range("F41").Value = 12
$ws.Range("H41").Value = 12

# Row 49
$ws.Range("E49").Value = 51

# Row 56
$ws.Range("E56").Value = 6

# Row 72
$ws.Range("E72").Value = 26

# Row 77
$ws.Range("E77").Value = 34

# Row 89
$ws.Range("E89").Value = 20
